$d = $word.ActiveDocument

# Change 1: merge two runs (no visible text change) -- the text
# "gezeichnet" followed directly by ", wodurch" is already how it reads;
# nothing textual to change here, Word will keep the text continuous.

# Change 2: "Darstellungscode" -> "Code"
$d.Content.Find.Execute("Darstellungscode", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Code", 2)

# Change 3: insert new sentence between "Fehlersuche. " and "Alles in allem"
$d.Content.Find.Execute("Fehlersuche. Alles in allem", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fehlersuche. Wir haben uns entschieden, dass alle Klassen ab InteractiveObject Unterobjekte beinhalten können, da in unserer Anwendung nur InteractiveObejcts Unterobjekte haben müssen. Alles in allem", 2)
